$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
}

# New weekly record: insert at row 11 (newest date), pushing existing rows down.
$ws.Rows(11).Insert()
Set-DataRow 11 8 "Terminal La Palmera de La Serena" "Coquimbo" 44426 4 100114007 "Jengibre" "Sin especificar" "Primera" 460 14000 15000 14500 "`$/caja 13 kilos" "Perú" 1115 13 "Hortaliza"

# Another new weekly record inserted further down the table (pushes rows below it down).
$ws.Rows(16).Insert()
Set-DataRow 16 8 "Terminal La Palmera de La Serena" "Coquimbo" 44419 4 100114007 "Jengibre" "Sin especificar" "Primera" 600 14000 15000 14500 "`$/caja 13 kilos" "Perú" 1115 13 "Hortaliza"

# New weekly record appended at the end of the table.
Set-DataRow 21 8 "Terminal La Palmera de La Serena" "Coquimbo" 44418 4 100114007 "Jengibre" "Sin especificar" "Primera" 500 14000 15000 14500 "`$/caja 13 kilos" "Perú" 1115 13 "Hortaliza"
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
